$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 335128992.9759778
$ws.Range("C2").Value = 413030972.2377607
$ws.Range("D2").Value = 490932951.49954426
$ws.Range("E2").Value = 568834930.7613283
$ws.Range("F2").Value = 646736910.0231119
$ws.Range("B3").Value = 730629184.4341357
$ws.Range("C3").Value = 808531163.6959186
$ws.Range("D3").Value = 886433142.9577022
$ws.Range("E3").Value = 964335122.2194864
$ws.Range("F3").Value = 1042237101.48127
$ws.Range("B4").Value = 1521982755.4793682
$ws.Range("C4").Value = 1599884734.7411513
$ws.Range("D4").Value = 1677786714.002935
$ws.Range("E4").Value = 1755688693.264719
$ws.Range("F4").Value = 1833590672.5265028
$ws.Range("B5").Value = 2472228971.695889
$ws.Range("C5").Value = 2550130950.957672
$ws.Range("D5").Value = 2628032930.2194557
$ws.Range("E5").Value = 2705934909.48124
$ws.Range("F5").Value = 2783836888.7430234
